$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 31.824752
$ws.Range("H2").Value = 95.47425600000001
$ws.Range("I2").Value = 0.886907633630525
$ws.Range("J2").Value = 0.886907633630525
$ws.Range("M2").Value = 15.50220733333333
$ws.Range("N2").Value = 46.506622
$ws.Range("O2").Value = 0.5994675913188158
$ws.Range("P2").Value = 0.5994675913188158
$ws.Range("Q2").Value = 493.3539038359147
$ws.Range("R2").Value = 4440.185134523233
$ws.Range("S2").Value = 0.5316723828547616
$ws.Range("T2").Value = 0.5316723828547616
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 31.824752
$ws.Range("H3").Value = 95.47425600000001
$ws.Range("I3").Value = 0.886907633630525
$ws.Range("J3").Value = 0.886907633630525
$ws.Range("O3").Value = 0.04399860030713892
$ws.Range("P3").Value = 0.04399860030713892
$ws.Range("Q3").Value = 36.21026647510401
$ws.Range("R3").Value = 325.892398275936
$ws.Range("S3").Value = 0.03902269448145988
$ws.Range("T3").Value = 0.03902269448145988
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 31.824752
$ws.Range("H4").Value = 95.47425600000001
$ws.Range("I4").Value = 0.886907633630525
$ws.Range("J4").Value = 0.886907633630525
$ws.Range("M4").Value = 8.848210666666667
$ws.Range("N4").Value = 26.544632
$ws.Range("O4").Value = 0.3421587275782868
$ws.Range("P4").Value = 0.3421587275782868
$ws.Range("Q4").Value = 281.5921101104213
$ws.Range("R4").Value = 2534.328990993792
$ws.Range("S4").Value = 0.3034631874024898
$ws.Range("T4").Value = 0.3034631874024898
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 31.824752
$ws.Range("H5").Value = 95.47425600000001
$ws.Range("I5").Value = 0.886907633630525
$ws.Range("J5").Value = 0.886907633630525
$ws.Range("M5").Value = 0.371739
$ws.Range("N5").Value = 1.115217
$ws.Range("O5").Value = 0.01437508079575842
$ws.Range("P5").Value = 0.01437508079575841
$ws.Range("Q5").Value = 11.830501483728
$ws.Range("R5").Value = 106.474513353552
$ws.Range("S5").Value = 0.0127493688918137
$ws.Range("T5").Value = 0.0127493688918137
$ws.Range("I6").Value = 0.06502043684278042
$ws.Range("J6").Value = 0.06502043684278042
$ws.Range("M6").Value = 15.50220733333333
$ws.Range("N6").Value = 46.506622
$ws.Range("O6").Value = 0.5994675913188158
$ws.Range("P6").Value = 0.5994675913188158
$ws.Range("Q6").Value = 36.16846346692467
$ws.Range("R6").Value = 325.516171202322
$ws.Range("S6").Value = 0.03897764466063876
$ws.Range("T6").Value = 0.03897764466063876
$ws.Range("I7").Value = 0.06502043684278042
$ws.Range("J7").Value = 0.06502043684278042
$ws.Range("O7").Value = 0.04399860030713892
$ws.Range("P7").Value = 0.04399860030713892
$ws.Range("S7").Value = 0.002860808212441065
$ws.Range("T7").Value = 0.002860808212441065
$ws.Range("I8").Value = 0.06502043684278042
$ws.Range("J8").Value = 0.06502043684278042
$ws.Range("M8").Value = 8.848210666666667
$ws.Range("N8").Value = 26.544632
$ws.Range("O8").Value = 0.3421587275782868
$ws.Range("P8").Value = 0.3421587275782868
$ws.Range("Q8").Value = 20.64391072598134
$ws.Range("R8").Value = 185.795196533832
$ws.Range("S8").Value = 0.02224730993671011
$ws.Range("T8").Value = 0.02224730993671011
$ws.Range("I9").Value = 0.06502043684278042
$ws.Range("J9").Value = 0.06502043684278042
$ws.Range("M9").Value = 0.371739
$ws.Range("N9").Value = 1.115217
$ws.Range("O9").Value = 0.01437508079575842
$ws.Range("P9").Value = 0.01437508079575841
$ws.Range("Q9").Value = 0.867310580463
$ws.Range("R9").Value = 7.805795224167
$ws.Range("S9").Value = 0.0009346740329904757
$ws.Range("T9").Value = 0.0009346740329904755
$ws.Range("G10").Value = 1.696588
$ws.Range("H10").Value = 5.089764000000001
$ws.Range("I10").Value = 0.04728133775640876
$ws.Range("J10").Value = 0.04728133775640876
$ws.Range("M10").Value = 15.50220733333333
$ws.Range("N10").Value = 46.506622
$ws.Range("O10").Value = 0.5994675913188158
$ws.Range("P10").Value = 0.5994675913188158
$ws.Range("Q10").Value = 26.30085893524534
$ws.Range("R10").Value = 236.707730417208
$ws.Range("S10").Value = 0.02834362965916574
$ws.Range("T10").Value = 0.02834362965916574
$ws.Range("G11").Value = 1.696588
$ws.Range("H11").Value = 5.089764000000001
$ws.Range("I11").Value = 0.04728133775640876
$ws.Range("J11").Value = 0.04728133775640876
$ws.Range("O11").Value = 0.04399860030713892
$ws.Range("P11").Value = 0.04399860030713892
$ws.Range("Q11").Value = 1.930381219576
$ws.Range("R11").Value = 17.373430976184
$ws.Range("S11").Value = 0.002080312681931066
$ws.Range("T11").Value = 0.002080312681931066
$ws.Range("G12").Value = 1.696588
$ws.Range("H12").Value = 5.089764000000001
$ws.Range("I12").Value = 0.04728133775640876
$ws.Range("J12").Value = 0.04728133775640876
$ws.Range("M12").Value = 8.848210666666667
$ws.Range("N12").Value = 26.544632
$ws.Range("O12").Value = 0.3421587275782868
$ws.Range("P12").Value = 0.3421587275782868
$ws.Range("Q12").Value = 15.01176803853867
$ws.Range("R12").Value = 135.105912346848
$ws.Range("S12").Value = 0.01617772236493203
$ws.Range("T12").Value = 0.01617772236493203
$ws.Range("G13").Value = 1.696588
$ws.Range("H13").Value = 5.089764000000001
$ws.Range("I13").Value = 0.04728133775640876
$ws.Range("J13").Value = 0.04728133775640876
$ws.Range("M13").Value = 0.371739
$ws.Range("N13").Value = 1.115217
$ws.Range("O13").Value = 0.01437508079575842
$ws.Range("P13").Value = 0.01437508079575841
$ws.Range("Q13").Value = 0.630687926532
$ws.Range("R13").Value = 5.676191338788001
$ws.Range("S13").Value = 0.0006796730503799189
$ws.Range("T13").Value = 0.0006796730503799188
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.02836866666666667
$ws.Range("H14").Value = 0.085106
$ws.Range("I14").Value = 0.0007905917702857979
$ws.Range("J14").Value = 0.0007905917702857978
$ws.Range("M14").Value = 15.50220733333333
$ws.Range("N14").Value = 46.506622
$ws.Range("O14").Value = 0.5994675913188158
$ws.Range("P14").Value = 0.5994675913188158
$ws.Range("Q14").Value = 0.4397769524368889
$ws.Range("R14").Value = 3.957992571932
$ws.Range("S14").Value = 0.0004739341442497058
$ws.Range("T14").Value = 0.0004739341442497057
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.02836866666666667
$ws.Range("H15").Value = 0.085106
$ws.Range("I15").Value = 0.0007905917702857979
$ws.Range("J15").Value = 0.0007905917702857978
$ws.Range("O15").Value = 0.04399860030713892
$ws.Range("P15").Value = 0.04399860030713892
$ws.Range("Q15").Value = 0.03227792567066667
$ws.Range("R15").Value = 0.290501331036
$ws.Range("S15").Value = 0.00003478493130691821
$ws.Range("T15").Value = 0.00003478493130691821
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.02836866666666667
$ws.Range("H16").Value = 0.085106
$ws.Range("I16").Value = 0.0007905917702857979
$ws.Range("J16").Value = 0.0007905917702857978
$ws.Range("M16").Value = 8.848210666666667
$ws.Range("N16").Value = 26.544632
$ws.Range("O16").Value = 0.3421587275782868
$ws.Range("P16").Value = 0.3421587275782868
$ws.Range("Q16").Value = 0.2510119389991111
$ws.Range("R16").Value = 2.259107450992
$ws.Range("S16").Value = 0.0002705078741548538
$ws.Range("T16").Value = 0.0002705078741548538
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.02836866666666667
$ws.Range("H17").Value = 0.085106
$ws.Range("I17").Value = 0.0007905917702857979
$ws.Range("J17").Value = 0.0007905917702857978
$ws.Range("M17").Value = 0.371739
$ws.Range("N17").Value = 1.115217
$ws.Range("O17").Value = 0.01437508079575842
$ws.Range("P17").Value = 0.01437508079575841
$ws.Range("Q17").Value = 0.010545739778
$ws.Range("R17").Value = 0.09491165800199999
$ws.Range("S17").Value = 0.00001136482057432002
$ws.Range("T17").Value = 0.00001136482057432002
